$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.49"
$ws.Range("E2").Value = "'-1.37%"
$ws.Range("E3").Value = "'-2.22%"
$ws.Range("D4").Value = "'4.879"
$ws.Range("E4").Value = "'1.45%"
$ws.Range("D5").Value = "'0.06323"
$ws.Range("E5").Value = "'-0.42%"
$ws.Range("D6").Value = "'6.922"
$ws.Range("E6").Value = "'-0.42%"
$ws.Range("D7").Value = "'1.275"
$ws.Range("E7").Value = "'33.09%"
$ws.Range("D8").Value = "'0.8741"
$ws.Range("E8").Value = "'-0.73%"
$ws.Range("D9").Value = "'0.1558"
$ws.Range("E9").Value = "'5.44%"
$ws.Range("D10").Value = "'0.05081"
$ws.Range("E10").Value = "'-1.26%"
$ws.Range("D11").Value = "'0.07505"
$ws.Range("E11").Value = "'2.81%"
$ws.Range("D12").Value = "'0.02968"
$ws.Range("E12").Value = "'-4.99%"
$ws.Range("D13").Value = "'0.09058"
$ws.Range("E13").Value = "'-0.18%"
$ws.Range("D14").Value = "'0.001583"
$ws.Range("E14").Value = "'1.41%"
$ws.Range("D15").Value = "'0.0006338"
$ws.Range("E15").Value = "'0.93%"
$ws.Range("D16").Value = "'0.005976"
$ws.Range("E16").Value = "'3.94%"
$ws.Range("D17").Value = "'3.455"
$ws.Range("E17").Value = "'0.03%"
$ws.Range("D18").Value = "'3.322"
$ws.Range("E18").Value = "'-2.27%"
$ws.Range("D19").Value = "'2.272"
$ws.Range("E19").Value = "'-1.01%"
$ws.Range("E20").Value = "'-0.13%"
$ws.Range("D21").Value = "'0.1335"
$ws.Range("E21").Value = "'3.18%"
$ws.Range("D22").Value = "'3.926"
$ws.Range("E22").Value = "'1.61%"
$ws.Range("D23").Value = "'0.04374"
$ws.Range("E23").Value = "'1.18%"
$ws.Range("D24").Value = "'0.001162"
$ws.Range("E24").Value = "'-1.55%"
$ws.Range("E25").Value = "'-1.85%"
$ws.Range("D26").Value = "'0.0001202"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("E27").Value = "'-4.39%"
$ws.Range("D40").Value = "'0.04110"
$ws.Range("E40").Value = "'0.32%"
$ws.Range("D41").Value = "'0.007043"
$ws.Range("E41").Value = "'5.83%"
$ws.Range("D42").Value = "'0.1172"
$ws.Range("E42").Value = "'0.58%"
$ws.Range("D43").Value = "'0.002254"
$ws.Range("E43").Value = "'2.38%"
$ws.Range("E44").Value = "'-12.75%"
$ws.Range("D45").Value = "'0.00005222"
$ws.Range("E45").Value = "'-0.10%"
$ws.Range("E47").Value = "'-11.19%"
